$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add a new row to the "Main" table for the new part ---
$lo = $ws.ListObjects.Item("Main")
$newListRow = $lo.ListRows.Add()
$r = $newListRow.Range.Row

# Copy the formatting of the previous data row onto the new row so that
# styles (date format, hyperlink style, currency format, etc.) match.
$ws.Range("A41:H41").Copy()
$ws.Range("A42:H42").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the new part's data ---
$ws.Cells.Item($r, 1).Value = "ZS-X10 V3 BLDC motor driver"
$ws.Cells.Item($r, 2).Value = 1
$ws.Cells.Item($r, 3).Value = "Ordered"
$ws.Cells.Item($r, 4).Value = 44.75
$ws.Cells.Item($r, 5).Formula = "=PRODUCT(B42*D42)"
$ws.Cells.Item($r, 6).Value = 45261
$ws.Cells.Item($r, 8).Value = "Allegro (kamami_pl)"

# Add the hyperlink for the LINK column.
$ws.Hyperlinks.Add($ws.Cells.Item($r, 7), "https://allegro.pl/oferta/sterownik-silnikow-bldc-6-20v-3a-13359462424") | Out-Null

# Re-apply the row formatting (Hyperlinks.Add resets the cell style); this
# restores the shared "LINK" column hyperlink style used throughout the table.
$ws.Range("A41:H41").Copy()
$ws.Range("A42:H42").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update the sheet view (scroll position / selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("H43").Select()

$excel.CalculateFullRebuild()

$wb.Save()
